$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 32, shifting existing rows 32:53 down to 33:54.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new data record.
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44505
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100101
$ws.Range("H32").Value = "Berries"
$ws.Range("I32").Value = 100101001
$ws.Range("J32").Value = "Arándano (blue)"
$ws.Range("K32").Value = "Sin especificar"
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 100
$ws.Range("N32").Value = 7000
$ws.Range("O32").Value = 7000
$ws.Range("P32").Value = 7000
$ws.Range("Q32").Value = "$/bandeja 2 kilos"
$ws.Range("R32").Value = "Provincia de Curicó"
$ws.Range("S32").Value = 3500
$ws.Range("T32").Value = 2
